# This script rolls the quarterly report forward by one quarter:
# - the 10 quarter-header labels in row 8 / row 24 shift to the next quarter
#   (oldest quarter "فصل دوم منتهی به 1399/06" drops off, "فصل چهارم منتهی به
#   1401/12" is appended)
# - every quarterly data row shifts its 10 values one column to the left and a
#   freshly-reported value is appended in column N
# - one data point (J14) is corrected in addition to the shift (read_price
#   algorithm change)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Roll the quarter headers forward (row 8 and row 24) ---
$headers = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(8, 5 + $i).Value = $headers[$i]
    $ws.Cells.Item(24, 5 + $i).Value = $headers[$i]
}

# --- Update quarterly data rows (E:N) ---
$row10 = @(2187021, 3572769, 1953376, 2389731, 899936, 2459757, 0, 2273286, 2286245, 1197237)
for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws.Cells.Item(10, 5 + $i).Value = $row10[$i]
}
$row11 = @(5935, -19043, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row11.Length; $i++) {
    $ws.Cells.Item(11, 5 + $i).Value = $row11[$i]
}
$row12 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, 5 + $i).Value = $row12[$i]
}
$row13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row13.Length; $i++) {
    $ws.Cells.Item(13, 5 + $i).Value = $row13[$i]
}
$row14 = @(4971, 8345, 2227, 6261, 15029, 81023, 9430, 49217, 29323, 84029)
for ($i = 0; $i -lt $row14.Length; $i++) {
    $ws.Cells.Item(14, 5 + $i).Value = $row14[$i]
}
$row15 = @(734, 456, 284, 572, 194, 1894, 475, 1140, 2157, -1455)
for ($i = 0; $i -lt $row15.Length; $i++) {
    $ws.Cells.Item(15, 5 + $i).Value = $row15[$i]
}
$row16 = @(10615, -6724, 2074, 11970, 7042, 7526, 10435, 10688, 16984, 8725)
for ($i = 0; $i -lt $row16.Length; $i++) {
    $ws.Cells.Item(16, 5 + $i).Value = $row16[$i]
}
$row17 = @(37110, 221964, 60885, 159732, 71315, 137564, 95464, 347705, 186479, 429337)
for ($i = 0; $i -lt $row17.Length; $i++) {
    $ws.Cells.Item(17, 5 + $i).Value = $row17[$i]
}
$row18 = @(-33, -40452, 0, 0, 0, 6191, 0, 0, 0, 505418)
for ($i = 0; $i -lt $row18.Length; $i++) {
    $ws.Cells.Item(18, 5 + $i).Value = $row18[$i]
}
$row19 = @(-38214, 358661, 249554, 256542, 470449, 506217, 1059884, -107878, 162667, 2056604)
for ($i = 0; $i -lt $row19.Length; $i++) {
    $ws.Cells.Item(19, 5 + $i).Value = $row19[$i]
}
$row20 = @(2208139, 4095976, 2268400, 2824808, 1463965, 3200172, 1175688, 2574158, 2683855, 4279895)
for ($i = 0; $i -lt $row20.Length; $i++) {
    $ws.Cells.Item(20, 5 + $i).Value = $row20[$i]
}
$row26 = @(978, 628, 960, 950, 950, 648, 648, 660, 660, 542)
for ($i = 0; $i -lt $row26.Length; $i++) {
    $ws.Cells.Item(26, 5 + $i).Value = $row26[$i]
}
$row27 = @(188, 532, 194, 195, 195, 695, 533, 543, 543, 877)
for ($i = 0; $i -lt $row27.Length; $i++) {
    $ws.Cells.Item(27, 5 + $i).Value = $row27[$i]
}
